# Update column F ("想去人数" / want-to-go count) values on the
# "展览", "演出" and "全部类型" sheets to match the refreshed data
# export (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1478
$ws1.Range("F3").Value  = 144
$ws1.Range("F4").Value  = 1765
$ws1.Range("F5").Value  = 33
$ws1.Range("F6").Value  = 146
$ws1.Range("F7").Value  = 659
$ws1.Range("F8").Value  = 35
$ws1.Range("F11").Value = 26
$ws1.Range("F13").Value = 150
$ws1.Range("F17").Value = 105
$ws1.Range("F18").Value = 4860
$ws1.Range("F19").Value = 48
$ws1.Range("F20").Value = 825
$ws1.Range("F21").Value = 108
$ws1.Range("F22").Value = 2221
$ws1.Range("F24").Value = 18
$ws1.Range("F25").Value = 2074

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 76

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1478
$ws4.Range("F3").Value  = 144
$ws4.Range("F4").Value  = 1765
$ws4.Range("F5").Value  = 33
$ws4.Range("F6").Value  = 146
$ws4.Range("F7").Value  = 659
$ws4.Range("F8").Value  = 35
$ws4.Range("F11").Value = 26
$ws4.Range("F13").Value = 150
$ws4.Range("F17").Value = 105
$ws4.Range("F18").Value = 4860
$ws4.Range("F19").Value = 76
$ws4.Range("F20").Value = 48
$ws4.Range("F22").Value = 825
$ws4.Range("F23").Value = 108
$ws4.Range("F24").Value = 2221
$ws4.Range("F26").Value = 18
$ws4.Range("F27").Value = 2074

$wb.Save()
